$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Relocate comments that are anchored below row 51, one row up, ---
# --- BEFORE deleting the row (comments do not auto-shift on Delete). ---
$cmt = $ws.Range("B53").Comment
$txt = $cmt.Text()
$cmt.Delete()
$ws.Range("B52").AddComment($txt) | Out-Null
$cmt = $ws.Range("C53").Comment
$txt = $cmt.Text()
$cmt.Delete()
$ws.Range("C52").AddComment($txt) | Out-Null
$cmt = $ws.Range("B71").Comment
$txt = $cmt.Text()
$cmt.Delete()
$ws.Range("B70").AddComment($txt) | Out-Null
$cmt = $ws.Range("B72").Comment
$txt = $cmt.Text()
$cmt.Delete()
$ws.Range("B71").AddComment($txt) | Out-Null
$cmt = $ws.Range("B74").Comment
$txt = $cmt.Text()
$cmt.Delete()
$ws.Range("B73").AddComment($txt) | Out-Null
$cmt = $ws.Range("B75").Comment
$txt = $cmt.Text()
$cmt.Delete()
$ws.Range("B74").AddComment($txt) | Out-Null
$cmt = $ws.Range("B76").Comment
$txt = $cmt.Text()
$cmt.Delete()
$ws.Range("B75").AddComment($txt) | Out-Null
$cmt = $ws.Range("B77").Comment
$txt = $cmt.Text()
$cmt.Delete()
$ws.Range("B76").AddComment($txt) | Out-Null
$cmt = $ws.Range("C77").Comment
$txt = $cmt.Text()
$cmt.Delete()
$ws.Range("C76").AddComment($txt) | Out-Null
$cmt = $ws.Range("B82").Comment
$txt = $cmt.Text()
$cmt.Delete()
$ws.Range("B81").AddComment($txt) | Out-Null
$cmt = $ws.Range("B83").Comment
$txt = $cmt.Text()
$cmt.Delete()
$ws.Range("B82").AddComment($txt) | Out-Null
$cmt = $ws.Range("B84").Comment
$txt = $cmt.Text()
$cmt.Delete()
$ws.Range("B83").AddComment($txt) | Out-Null
$cmt = $ws.Range("B85").Comment
$txt = $cmt.Text()
$cmt.Delete()
$ws.Range("B84").AddComment($txt) | Out-Null
$cmt = $ws.Range("B86").Comment
$txt = $cmt.Text()
$cmt.Delete()
$ws.Range("B85").AddComment($txt) | Out-Null
$cmt = $ws.Range("C86").Comment
$txt = $cmt.Text()
$cmt.Delete()
$ws.Range("C85").AddComment($txt) | Out-Null
$cmt = $ws.Range("B101").Comment
$txt = $cmt.Text()
$cmt.Delete()
$ws.Range("B100").AddComment($txt) | Out-Null
$cmt = $ws.Range("B103").Comment
$txt = $cmt.Text()
$cmt.Delete()
$ws.Range("B102").AddComment($txt) | Out-Null
$cmt = $ws.Range("B104").Comment
$txt = $cmt.Text()
$cmt.Delete()
$ws.Range("B103").AddComment($txt) | Out-Null
$cmt = $ws.Range("B105").Comment
$txt = $cmt.Text()
$cmt.Delete()
$ws.Range("B104").AddComment($txt) | Out-Null
$cmt = $ws.Range("B106").Comment
$txt = $cmt.Text()
$cmt.Delete()
$ws.Range("B105").AddComment($txt) | Out-Null
$cmt = $ws.Range("B107").Comment
$txt = $cmt.Text()
$cmt.Delete()
$ws.Range("B106").AddComment($txt) | Out-Null
$cmt = $ws.Range("B108").Comment
$txt = $cmt.Text()
$cmt.Delete()
$ws.Range("B107").AddComment($txt) | Out-Null

# --- Update the M07 comment (continuous delivery pipeline): drop SBoM
# --- line (renumber), drop the two ISD/Jenkins paragraphs. ---
$m07 = @'
M07: Het project gebruikt een continuous delivery pipeline om het product te bouwen, testen en op te leveren

Er is een geautomatiseerde continuous delivery pipeline die aantoonbaar correct werkt en de software bouwt, installeert in de testomgevingen, test op functionele en niet-functionele eigenschappen en oplevert, al dan niet inclusief installatie in de productieomgeving.

De geautomatiseerde continuous delivery pipeline voert ten minste de volgende activiteiten uit:

1. Bouw van de software,
2. Unit tests,
3. Regressietests,
4. Beveiligingstests,
5. Performancetests,
6. Toegankelijkheidstests,
7. Broncodekwaliteitscontroles,
8. Installatie van de software in test, acceptatie en/of productieomgevingen,
9. Oplevering van het totale product, dus inclusief alle deliverables, in de vorm zoals bruikbaar voor en afgesproken met de opdrachtgevende organisatie.

Performance- en beveiligingstests op de software zijn ook onderdeel van de continuous delivery pipeline, maar vanwege doorlooptijden en licenties is dat niet altijd haalbaar; in dat geval vinden de performance- en beveiligingstests zo veel mogelijk, en bij voorkeur dagelijks, plaats. Performance- en beveiligingstests op de software vinden plaats in de testomgeving van het project. Als ICTU verantwoordelijk is voor het operationeel beheer laat ICTU de performance- en beveiligingstesten op de software (ook) uitvoeren in een productie-like omgeving.

Niet alle testen en controles kunnen altijd geautomatiseerd worden uitgevoerd. Denk aan kwaliteitscontroles op architectuurbeslissingen of het testen van toegankelijkheidseisen. Waar mogelijk wordt wel een zo groot mogelijk deel van de testen en controles geautomatiseerd en als onderdeel van de pipeline uitgevoerd.

Rationale

Software incrementeel opleveren vereist dat de software frequent gebouwd, getest en opgeleverd kan worden. Om dit efficiënt en foutvrij te doen, dient het proces van bouwen, testen en opleveren geautomatiseerd te zijn; een continuous delivery pipeline faciliteert dit.


'@
$ws.Range("B42").Comment.Text($m07)

# --- Delete row 51 ("9. Produceren van een SBoM" line); everything below
# --- shifts up one row. ---
$ws.Rows("51").Delete()

# --- Fix the renumbered line that is now row 51 (was row 52). ---
$ws.Range("B51").Value2 = "9. Oplevering van het totale product, dus inclusief alle deliverables, in de vorm zoals bruikbaar voor en afgesproken met de opdrachtgevende organisatie."

# --- M16 task list (rows 52..68): drop the tool-name suffix, capitalise. ---
$ws.Range("B52").Value2 = "1. Product en sprint backlog management en agile werken"
$ws.Range("B53").Value2 = "2. Inrichten en uitvoeren van een continuous delivery pipeline"
$ws.Range("B54").Value2 = "3. Monitoren van de kwaliteit van broncode"
$ws.Range("B55").Value2 = "4. Versiebeheer van op te leveren producten"
$ws.Range("B56").Value2 = "5. Release van software"
$ws.Range("B57").Value2 = "6. Maken van testrapportages"
$ws.Range("B58").Value2 = "7. Maken van kwaliteitsrapportages"
$ws.Range("B59").Value2 = "8. Controleren op aanwezigheid van bekende kwetsbaarheden in externe software"
$ws.Range("B60").Value2 = "9. Statische controle van de software op aanwezigheid van kwetsbare constructies"
$ws.Range("B61").Value2 = "10. Dynamische controle van de software op aanwezigheid van kwetsbare constructies"
$ws.Range("B62").Value2 = "11. Controleren van container images op aanwezigheid van bekende kwetsbaarheden"
$ws.Range("B63").Value2 = "12. Testen van performance en schaalbaarheid"
$ws.Range("B64").Value2 = "13. Testen op toegankelijkheid van de applicatie"
$ws.Range("B65").Value2 = "14. Produceren van een `"software bill of materials`" (SBoM)"
$ws.Range("B66").Value2 = "15. Opslaan van artifacten"
$ws.Range("B67").Value2 = "16. Registratie van incidenten bij gebruik en beheer"
$ws.Range("B68").Value2 = "17. Bij het uitvoeren van operationeel beheer; uitrollen van de software in de productieomgeving"

# --- Rewrite the M16 comment (now anchored at B52) with the merged table. ---
$m16 = @'
M16: Het project gebruikt tools voor vastgestelde taken

Voor vastgestelde taken bij het ontwikkelen, onderhouden en operationeel beheren van software, stelt ICTU het gebruik van tools verplicht. ICTU adviseert per taak specifieke tools en ondersteunt projecten bij het gebruik daarvan.

ICTU adviseert en ondersteunt voor de hieronder genoemde taken specifieke tools. Projecten gebruiken deze tools, of gelijkwaardige alternatieven.

Activiteit                                                                                   Tools                                                                                    
Product en sprint backlog management en agile werken                                         Azure DevOps of Jira                                                                     
Inrichten en uitvoeren van een continuous delivery pipeline                                  Jenkins, GitLab CI/CD (Continuous Integration, Delivery, and Deployment) of Azure DevOps 
Monitoren van de kwaliteit van broncode                                                      SonarQube                                                                                
Versiebeheer van op te leveren producten                                                     GitLab of Azure DevOps                                                                   
Release van software                                                                         Releaseserver in het ontwikkelplatform                                                   
Maken van testrapportages                                                                    JUnit, Robot Framework, TestNG, of hiermee compatible tools                              
Maken van kwaliteitsrapportages                                                              Quality-time                                                                             
Controleren op aanwezigheid van bekende kwetsbaarheden in externe software                   OWASP (Open Web Application Security Project) Dependency-Check en/of Dependency-Track    
Statische controle van de software op aanwezigheid van kwetsbare constructies                SonarQube                                                                                
Dynamische controle van de software op aanwezigheid van kwetsbare constructies               ZAP (Zed Attack Proxy) by Checkmarx                                                      
Controleren van container images op aanwezigheid van bekende kwetsbaarheden                  Trivy                                                                                    
Testen van performance en schaalbaarheid                                                     JMeter en Performancetestrunner                                                          
Testen op toegankelijkheid van de applicatie                                                 Axe                                                                                      
Produceren van een "software bill of materials" (SBoM)                                       Tools die een SBoM in CycloneDX-formaat (zie https://cyclonedx.org) genereren            
Opslaan van artifacten                                                                       Nexus of Harbor                                                                          
Registratie van incidenten bij gebruik en beheer                                             Jira                                                                                     
Bij het uitvoeren van operationeel beheer; uitrollen van de software in de productieomgeving Ansible                                                                                  

N.B. Onder het ondersteunen van "agile werken" vallen het opvoeren van eisen, het opvoeren van logische testgevallen, het koppelen van logische testgevallen aan eisen, het bijhouden van een werkvoorraad, het plannen van iteraties en het toewijzen van eisen aan iteraties. De 'eisen' worden, conform Scrumterminologie, geregistreerd als epics en/of user stories, de werkvoorraad als product backlog en de iteraties als sprints. Het toewijzen van eisen aan iteraties gebeurt via de sprint backlog.

Rationale

Projecten hebben een redelijke vrijheid bij het kiezen en gebruiken van tools, maar voor een aantal taken is het gebruik verplicht gesteld. Deze tools zijn nodig voor een efficiënte uitvoering van de Kwaliteitsaanpak. Uniform gebruik van deze tools maakt het mogelijk koppeling tussen die tools voor alle projecten te standaardiseren; daarnaast bevordert het de uitwisselbaarheid van medewerkers en neemt het risico op het gebruik van onvolwassen tools af. Tot slot is het gebruik in een aantal gevallen, ten behoeve van informatiebeveiliging bij de overheid, verplicht.


'@
$ws.Range("B52").Comment.Text($m16)

# --- Rewrite the M18 comment (now anchored at B104): 'bij M16' -> 'in M16',
# --- plus appended sentence about project responsibility. ---
$m18 = @'
M18: ICTU biedt ondersteuning voor verplicht gestelde tools

ICTU zorgt voor technische en functionele ondersteuning aan projecten bij het gebruik van alle verplichte tools.

ICTU zorgt voor ondersteuning van de in M16: Het project gebruikt tools voor vastgestelde taken verplicht gestelde tools. Een team van specialisten met kennis, ervaring en capaciteit is beschikbaar voor ondersteuning aan projecten. Projecten zijn verantwoordelijk voor de correcte werking van de pipeline.

Bij de selectie van tools ter ondersteuning van de projectuitvoering geeft ICTU de voorkeur aan open source tools. Ook tools die ICTU zelf ontwikkelt ter ondersteuning van softwareontwikkelprojecten worden bij voorkeur open source beschikbaar gesteld.

Rationale

De keuze om het gebruik van een aantal tools verplicht te stellen (M16: Het project gebruikt tools voor vastgestelde taken) volgt uit de belangrijke rol die die tools spelen in de ontwikkelstraat en in Quality-time, het kwaliteitssysteem van ICTU. Met de verplichting komt ook een verantwoordelijkheid: om projecten in staat te stellen snel en effectief met deze tools te werken, moeten die projecten ondersteund worden.

De verplicht gestelde tools zijn beperkt in aantal, bewezen en gangbaar; veel medewerkers zullen deze tools al kennen.

De voorkeur voor open source tools is conform de rationale uit NORA (Nederlandse Overheid Referentiearchitectuur) voor het gebruik van open source tools, zoals beschreven in NORA v3.0 drijfveer "Beleid open standaarden". De voorkeur voor het open source beschikbaar stellen van eigen ontwikkelde tools is conform de "Beleidsbrief vrijgeven van de broncode van overheidssoftware" van de staatssecretaris van Binnenlandse Zaken en Koninkrijksrelaties, 17 april 2020.


'@
$ws.Range("B104").Comment.Text($m18)
